# Updated symbol list on Sat Jan 14 04:39:15 UTC 2023 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures for
# the coin rows on the active worksheet. The source values are numeric-
# looking strings (e.g. "306.68", "5.96%") that must stay stored as plain
# text, exactly like the original cells, instead of being auto-converted
# to numbers/percentages by Excel. Forcing the cell to Text format ("@")
# before writing the value preserves that literal text representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "306.68"
    "E2"  = "5.96%"
    "E3"  = "9.06%"
    "D4"  = "5.327"
    "E4"  = "4.04%"
    "D5"  = "0.07459"
    "D6"  = "7.750"
    "E6"  = "5.42%"
    "D7"  = "3.703"
    "E7"  = "8.75%"
    "D8"  = "1.556"
    "E8"  = "15.29%"
    "D9"  = "0.9228"
    "E9"  = "0.69%"
    "D10" = "0.01649"
    "E10" = "2,453.14%"
    "D11" = "0.1665"
    "E11" = "4.78%"
    "D12" = "0.07579"
    "E12" = "13.45%"
    "D13" = "0.07959"
    "E13" = "3.57%"
    "D14" = "0.03080"
    "E14" = "4.73%"
    "D15" = "0.09839"
    "E15" = "9.42%"
    "E16" = "-3.37%"
    "D17" = "0.04548"
    "D18" = "0.006429"
    "E18" = "2.37%"
    "D19" = "3.476"
    "E19" = "0.66%"
    "D20" = "2.242"
    "E20" = "0.95%"
    "D21" = "0.3277"
    "E21" = "1.99%"
    "D22" = "0.1326"
    "E22" = "1.27%"
    "D23" = "4.209"
    "E23" = "3.53%"
    "D24" = "0.1629"
    "E24" = "5.06%"
    "D25" = "0.001194"
    "E25" = "0.41%"
    "D26" = "0.004533"
    "D27" = "0.0001168"
    "E27" = "-6.47%"
    "E28" = "2.93%"
    "D40" = "0.04518"
    "E40" = "6.78%"
    "D41" = "0.007413"
    "E41" = "9.83%"
    "D42" = "0.1365"
    "E42" = "9.98%"
    "D43" = "0.002256"
    "E43" = "14.05%"
    "D44" = "0.01382"
    "E44" = "17.79%"
    "D45" = "0.00006006"
    "E45" = "7.02%"
    "D46" = "1.892"
    "E46" = "-3.82%"
    "D47" = "0.01298"
    "E47" = "-0.61%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Keep the new value as literal text (matching the original inline
    # string cells) rather than letting Excel coerce it into a number or
    # a percentage.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
